$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the table text: "Math derived from Reactions" -> "Implied math for rates of change of species"
$ws.Range("A14").Value = "Implied math for rates of change of species"

# Move the selection/active cell to A17 (was C24 / A1:C24)
[void]$ws.Range("A17").Select()
